# Apply the edit described by the diff:
# - The 6 data rows (rows 2-7) are replaced by a reordered/updated set of
#   5 data rows (rows 2-6); the former row 7 is removed, shrinking the
#   sheet from A1:F7 to A1:F6.
# - Column D ("Advogado") width changes from 59 to 48 characters.
# - Some "OAB" (column E) values flip between numeric and text storage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2..6, in column order A..F.
# E.IsText marks OAB values that must be stored as text instead of a number.
# F.Text marks "date-looking" strings that must stay literal text (and not
# get auto-converted to a real date) by using Excel's leading-apostrophe,
# force-text entry convention.
$data = @(
    @{ A = "processo_3";      B = "3130687-11.2024.8.01.5042"; C = "Nome Aleatório 98"; D = "Advogado Exemplo ";  E = 43679;   EIsText = $false; F = "7/8/2024" },
    @{ A = "processo_oab";    B = "5056353-90.2024.8.12.0001"; C = "João da Silva";      D = "Maria Souza`nData de Distribuição: 18/04/2024`n"; E = $null; EIsText = $false; F = "18/04/2024" },
    @{ A = "processo_adv";    B = "1855099-63.2024.8.12.0001"; C = "João da Silva";      D = $null;             E = 45950;   EIsText = $false; F = "16/08/2024" },
    @{ A = "processo_1";      B = "3781128-20.2024.8.01.8252"; C = "Nome Aleatório 86"; D = "Advogado Exemplo ";  E = 44432;   EIsText = $false; F = "12/5/2024" },
    @{ A = "copy_processo_2"; B = "4835245-15.2024.8.01.2832"; C = "Nome Aleatório 2";  D = "Advogado Exemplo ";  E = "12723"; EIsText = $true;  F = "25/5/2024" }
)

$rowIndex = 2
foreach ($rec in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $rec.A
    $ws.Cells.Item($rowIndex, 2).Value = $rec.B
    $ws.Cells.Item($rowIndex, 3).Value = $rec.C

    if ($null -eq $rec.D) {
        $ws.Cells.Item($rowIndex, 4).Value = $null
    } else {
        $ws.Cells.Item($rowIndex, 4).Value = $rec.D
    }

    if ($null -eq $rec.E) {
        $ws.Cells.Item($rowIndex, 5).Value = $null
    } elseif ($rec.EIsText) {
        $ws.Cells.Item($rowIndex, 5).Value = "'" + $rec.E
    } else {
        $ws.Cells.Item($rowIndex, 5).Value = [double]$rec.E
    }

    # Force-text entry (leading apostrophe) so day/month-like strings are
    # kept as literal text instead of being parsed into date serials.
    $ws.Cells.Item($rowIndex, 6).Value = "'" + $rec.F

    # Multi-line text (embedded `n) can trigger an automatic custom row
    # height; AutoFit() re-settles the row back to the sheet's normal
    # (non-custom) height so we don't leave stray row-height overrides.
    $ws.Rows.Item($rowIndex).AutoFit()

    $rowIndex = $rowIndex + 1
}

# Remove the now-unused trailing row (previously row 7).
$ws.Rows.Item(7).Delete()

# Column D width change (59 -> 48 characters). Excel's ColumnWidth property
# round-trips with a +5/6 character offset versus the raw stored XML width,
# so compensate to land exactly on 48.
$ws.Columns.Item(4).ColumnWidth = 48 - (5 / 6)
